$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for new column K
$ws.Range("K1").Value = "REGION_CENSUS"

# Seed shared-string table in the same order as the source region lookup table
$ws.Cells.Item(2, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(5, 11).Value = "Superior CA"
$ws.Cells.Item(9, 11).Value = "North Coast"
$ws.Cells.Item(3, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(28, 11).Value = "Central Coast"
$ws.Cells.Item(11, 11).Value = "Southern San Joaquin Valley"
$ws.Cells.Item(34, 11).Value = "Inland Empire"
$ws.Cells.Item(20, 11).Value = "Los Angeles County"
$ws.Cells.Item(31, 11).Value = "Orange County"
$ws.Cells.Item(14, 11).Value = "San Diego - Imperial"

# Remaining rows (2020 multi-year census region groupings)
$ws.Cells.Item(4, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(6, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(7, 11).Value = "Superior CA"
$ws.Cells.Item(8, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(10, 11).Value = "Superior CA"
$ws.Cells.Item(12, 11).Value = "Superior CA"
$ws.Cells.Item(13, 11).Value = "North Coast"
$ws.Cells.Item(15, 11).Value = "Southern San Joaquin Valley"
$ws.Cells.Item(16, 11).Value = "Southern San Joaquin Valley"
$ws.Cells.Item(17, 11).Value = "Southern San Joaquin Valley"
$ws.Cells.Item(18, 11).Value = "North Coast"
$ws.Cells.Item(19, 11).Value = "Superior CA"
$ws.Cells.Item(21, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(22, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(23, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(24, 11).Value = "North Coast"
$ws.Cells.Item(25, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(26, 11).Value = "Superior CA"
$ws.Cells.Item(27, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(29, 11).Value = "North Coast"
$ws.Cells.Item(30, 11).Value = "Superior CA"
$ws.Cells.Item(32, 11).Value = "Superior CA"
$ws.Cells.Item(33, 11).Value = "Superior CA"
$ws.Cells.Item(35, 11).Value = "Superior CA"
$ws.Cells.Item(36, 11).Value = "Central Coast"
$ws.Cells.Item(37, 11).Value = "Inland Empire"
$ws.Cells.Item(38, 11).Value = "San Diego - Imperial"
$ws.Cells.Item(39, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(40, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(41, 11).Value = "Central Coast"
$ws.Cells.Item(42, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(43, 11).Value = "Central Coast"
$ws.Cells.Item(44, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(45, 11).Value = "Central Coast"
$ws.Cells.Item(46, 11).Value = "Superior CA"
$ws.Cells.Item(47, 11).Value = "Superior CA"
$ws.Cells.Item(48, 11).Value = "Superior CA"
$ws.Cells.Item(49, 11).Value = "San Francisco Bay Area"
$ws.Cells.Item(50, 11).Value = "North Coast"
$ws.Cells.Item(51, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(52, 11).Value = "Superior CA"
$ws.Cells.Item(53, 11).Value = "Superior CA"
$ws.Cells.Item(54, 11).Value = "North Coast"
$ws.Cells.Item(55, 11).Value = "Southern San Joaquin Valley"
$ws.Cells.Item(56, 11).Value = "Northern San Joaquin Valley"
$ws.Cells.Item(57, 11).Value = "Central Coast"
$ws.Cells.Item(58, 11).Value = "Superior CA"
$ws.Cells.Item(59, 11).Value = "Superior CA"

# Match column width used for the other lookup columns
$ws.Columns.Item(11).ColumnWidth = 32.333333333333336

# Leave the selection where the author left it when saving
$ws.Range("N51").Select() | Out-Null
